$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 37.63904266666666
$ws.Range("H2").Value = 112.917128
$ws.Range("I2").Value = 0.4850220755088102
$ws.Range("J2").Value = 0.4850220755088102
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 98.946724
$ws.Range("N2").Value = 296.840172
$ws.Range("O2").Value = 0.2098009692989996
$ws.Range("P2").Value = 0.2098009692989996
$ws.Range("Q2").Value = 3724.25996636289
$ws.Range("R2").Value = 33518.33969726601
$ws.Range("S2").Value = 0.1017581015731609
$ws.Range("T2").Value = 0.101758101573161
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 37.63904266666666
$ws.Range("H3").Value = 112.917128
$ws.Range("I3").Value = 0.4850220755088102
$ws.Range("J3").Value = 0.4850220755088102
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 163.0062356666667
$ws.Range("N3").Value = 489.018707
$ws.Range("O3").Value = 0.345629090707923
$ws.Range("P3").Value = 0.3456290907079231
$ws.Range("Q3").Value = 6135.398659190388
$ws.Range("R3").Value = 55218.58793271349
$ws.Range("S3").Value = 0.1676377389313796
$ws.Range("T3").Value = 0.1676377389313797
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 37.63904266666666
$ws.Range("H4").Value = 112.917128
$ws.Range("I4").Value = 0.4850220755088102
$ws.Range("J4").Value = 0.4850220755088102
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 65.39610666666668
$ws.Range("N4").Value = 196.18832
$ws.Range("O4").Value = 0.1386621609326595
$ws.Range("P4").Value = 0.1386621609326595
$ws.Range("Q4").Value = 2461.446849060551
$ws.Range("R4").Value = 22153.02164154496
$ws.Range("S4").Value = 0.06725420909009514
$ws.Range("T4").Value = 0.06725420909009516
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 37.63904266666666
$ws.Range("H5").Value = 112.917128
$ws.Range("I5").Value = 0.4850220755088102
$ws.Range("J5").Value = 0.4850220755088102
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 144.2727966666667
$ws.Range("N5").Value = 432.81839
$ws.Range("O5").Value = 0.3059077790604178
$ws.Range("P5").Value = 0.3059077790604179
$ws.Range("Q5").Value = 5430.289949375991
$ws.Range("R5").Value = 48872.60954438392
$ws.Range("S5").Value = 0.1483720259141744
$ws.Range("T5").Value = 0.1483720259141744
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.57434666666667
$ws.Range("H6").Value = 52.72304
$ws.Range("I6").Value = 0.2264655392929762
$ws.Range("J6").Value = 0.2264655392929762
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 98.946724
$ws.Range("N6").Value = 296.840172
$ws.Range("O6").Value = 0.2098009692989996
$ws.Range("P6").Value = 0.2098009692989996
$ws.Range("Q6").Value = 1738.924029106987
$ws.Range("R6").Value = 15650.31626196288
$ws.Range("S6").Value = 0.04751268965648708
$ws.Range("T6").Value = 0.04751268965648708
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.57434666666667
$ws.Range("H7").Value = 52.72304
$ws.Range("I7").Value = 0.2264655392929762
$ws.Range("J7").Value = 0.2264655392929762
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.0062356666667
$ws.Range("N7").Value = 489.018707
$ws.Range("O7").Value = 0.345629090707923
$ws.Range("P7").Value = 0.3456290907079231
$ws.Range("Q7").Value = 2864.728094434365
$ws.Range("R7").Value = 25782.55284990928
$ws.Range("S7").Value = 0.07827307842251077
$ws.Range("T7").Value = 0.07827307842251079
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 17.57434666666667
$ws.Range("H8").Value = 52.72304
$ws.Range("I8").Value = 0.2264655392929762
$ws.Range("J8").Value = 0.2264655392929762
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 65.39610666666668
$ws.Range("N8").Value = 196.18832
$ws.Range("O8").Value = 0.1386621609326595
$ws.Range("P8").Value = 0.1386621609326595
$ws.Range("Q8").Value = 1149.293849210311
$ws.Range("R8").Value = 10343.6446428928
$ws.Range("S8").Value = 0.03140220105514418
$ws.Range("T8").Value = 0.03140220105514419
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 17.57434666666667
$ws.Range("H9").Value = 52.72304
$ws.Range("I9").Value = 0.2264655392929762
$ws.Range("J9").Value = 0.2264655392929762
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 144.2727966666667
$ws.Range("N9").Value = 432.81839
$ws.Range("O9").Value = 0.3059077790604178
$ws.Range("P9").Value = 0.3059077790604179
$ws.Range("Q9").Value = 2535.500143189511
$ws.Range("R9").Value = 22819.5012887056
$ws.Range("S9").Value = 0.06927757015883414
$ws.Range("T9").Value = 0.06927757015883415
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 19.168158
$ws.Range("H10").Value = 57.504474
$ws.Range("I10").Value = 0.2470036195972184
$ws.Range("J10").Value = 0.2470036195972184
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 98.946724
$ws.Range("N10").Value = 296.840172
$ws.Range("O10").Value = 0.2098009692989996
$ws.Range("P10").Value = 0.2098009692989996
$ws.Range("Q10").Value = 1896.626439214392
$ws.Range("R10").Value = 17069.63795292953
$ws.Range("S10").Value = 0.05182159881185779
$ws.Range("T10").Value = 0.05182159881185779
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 19.168158
$ws.Range("H11").Value = 57.504474
$ws.Range("I11").Value = 0.2470036195972184
$ws.Range("J11").Value = 0.2470036195972184
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 163.0062356666667
$ws.Range("N11").Value = 489.018707
$ws.Range("O11").Value = 0.345629090707923
$ws.Range("P11").Value = 0.3456290907079231
$ws.Range("Q11").Value = 3124.529280243902
$ws.Range("R11").Value = 28120.76352219512
$ws.Range("S11").Value = 0.08537163644295231
$ws.Range("T11").Value = 0.08537163644295233
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 19.168158
$ws.Range("H12").Value = 57.504474
$ws.Range("I12").Value = 0.2470036195972184
$ws.Range("J12").Value = 0.2470036195972184
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 65.39610666666668
$ws.Range("N12").Value = 196.18832
$ws.Range("O12").Value = 0.1386621609326595
$ws.Range("P12").Value = 0.1386621609326595
$ws.Range("Q12").Value = 1253.522905171521
$ws.Range("R12").Value = 11281.70614654368
$ws.Range("S12").Value = 0.03425005565153889
$ws.Range("T12").Value = 0.0342500556515389
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 19.168158
$ws.Range("H13").Value = 57.504474
$ws.Range("I13").Value = 0.2470036195972184
$ws.Range("J13").Value = 0.2470036195972184
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 144.2727966666667
$ws.Range("N13").Value = 432.81839
$ws.Range("O13").Value = 0.3059077790604178
$ws.Range("P13").Value = 0.3059077790604179
$ws.Range("Q13").Value = 2765.44376160854
$ws.Range("R13").Value = 24888.99385447686
$ws.Range("S13").Value = 0.07556032869086937
$ws.Range("T13").Value = 0.07556032869086939
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.221194
$ws.Range("H14").Value = 9.663582
$ws.Range("I14").Value = 0.04150876560099527
$ws.Range("J14").Value = 0.04150876560099527
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 98.946724
$ws.Range("N14").Value = 296.840172
$ws.Range("O14").Value = 0.2098009692989996
$ws.Range("P14").Value = 0.2098009692989996
$ws.Range("Q14").Value = 318.726593668456
$ws.Range("R14").Value = 2868.539343016104
$ws.Range("S14").Value = 0.008708579257493779
$ws.Range("T14").Value = 0.008708579257493779
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.221194
$ws.Range("H15").Value = 9.663582
$ws.Range("I15").Value = 0.04150876560099527
$ws.Range("J15").Value = 0.04150876560099527
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 163.0062356666667
$ws.Range("N15").Value = 489.018707
$ws.Range("O15").Value = 0.345629090707923
$ws.Range("P15").Value = 0.3456290907079231
$ws.Range("Q15").Value = 525.0747082920527
$ws.Range("R15").Value = 4725.672374628474
$ws.Range("S15").Value = 0.01434663691108031
$ws.Range("T15").Value = 0.01434663691108031
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.221194
$ws.Range("H16").Value = 9.663582
$ws.Range("I16").Value = 0.04150876560099527
$ws.Range("J16").Value = 0.04150876560099527
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 65.39610666666668
$ws.Range("N16").Value = 196.18832
$ws.Range("O16").Value = 0.1386621609326595
$ws.Range("P16").Value = 0.1386621609326595
$ws.Range("Q16").Value = 210.6535464180267
$ws.Range("R16").Value = 1895.88191776224
$ws.Range("S16").Value = 0.005755695135881245
$ws.Range("T16").Value = 0.005755695135881246
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.221194
$ws.Range("H17").Value = 9.663582
$ws.Range("I17").Value = 0.04150876560099527
$ws.Range("J17").Value = 0.04150876560099527
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 144.2727966666667
$ws.Range("N17").Value = 432.81839
$ws.Range("O17").Value = 0.3059077790604178
$ws.Range("P17").Value = 0.3059077790604179
$ws.Range("Q17").Value = 464.7306669858867
$ws.Range("R17").Value = 4182.57600287298
$ws.Range("S17").Value = 0.01269785429653993
$ws.Range("T17").Value = 0.01269785429653993
